# Clear out the example input/output data rows (chemical names, numeric
# values and units) that were used as sample data, keeping the header row
# (Input/output | In/out | Units) intact. This also prunes the now-unused
# shared strings (chemical names, "kg", "kg ", "Electricity", "kWh").
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2:C48").ClearContents()

# Extend the "Output" cell style to column C (rows 1-39) so the formerly
# empty/unstyled Units column matches the styled look of columns A and B.
$ws.Range("C1:C39").Style = "Output"

# Match the author's new selection left after clearing the example rows.
$ws.Range("A2:C8").Select() | Out-Null
